$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# Each cell is temporarily switched to Text number format before the
# value is written (so digit/dot-heavy strings such as "591.06" or
# "0.601" are not auto-converted to numeric values - they are plain
# text in the source data, same as every other cell in this column),
# then the style is reset back to "Normal" so no new/visible cell
# formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.578.88'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.496.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.06'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.492.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.103.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.582.03'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.490.90'
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '392.87'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.97'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.537'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000122'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.24'
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.47'
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.65'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.38'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.882'
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '27.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0747'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '26.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.795.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '350.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.67'
$ws.Range("D51").Style = "Normal"

# --- Column B/C updates (RenderToken / InjectiveProtocol swap rows 40-41) ---
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'

# --- Column E (Volume(1h) %) updates ---
$ws.Range("E2").Value = '  +3.83%  '
$ws.Range("E3").Value = '  +1.74%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +3.24%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +8.50%  '
$ws.Range("E9").Value = '  +1.59%  '
$ws.Range("E10").Value = '  +6.71%  '
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("E12").Value = '  +3.44%  '
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  +3.68%  '
$ws.Range("E16").Value = '  +2.35%  '
$ws.Range("E17").Value = '  +3.70%  '
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("E19").Value = '  +2.88%  '
$ws.Range("E20").Value = '  +3.35%  '
$ws.Range("E21").Value = '  +3.73%  '
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E25").Value = '  +3.42%  '
$ws.Range("E26").Value = '  +4.86%  '
$ws.Range("E27").Value = '  +7.05%  '
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("E30").Value = '  +3.82%  '
$ws.Range("E31").Value = '  +4.63%  '
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("E33").Value = '  +2.67%  '
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +7.68%  '
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E38").Value = '  +2.26%  '
$ws.Range("E39").Value = '  +4.78%  '
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("E41").Value = '  +4.02%  '
$ws.Range("E42").Value = '  +5.50%  '
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("E44").Value = '  +1.88%  '
$ws.Range("E45").Value = '  -0.53%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("E48").Value = '  +2.84%  '
$ws.Range("E49").Value = '  +5.52%  '
$ws.Range("E50").Value = '  +4.92%  '
$ws.Range("E51").Value = '  +12.02%  '
